$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2333.95
$ws.Range("I40").Value = 3079.6
$ws.Range("J40").Value = 1588.3
$ws.Range("K40").Value = 3079.6
$ws.Range("L40").Value = 1588.3
$ws.Range("M40").Value = -2904.6
$ws.Range("N40").Value = -1938.3
$ws.Range("H54").Value = 3076
$ws.Range("I54").Value = 3076
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 3076
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0
$ws.Range("M54").Value = -2590
$ws.Range("H64").Value = 3859.25
$ws.Range("I64").Value = 3441.7856
$ws.Range("J64").Value = 4833.3335
$ws.Range("K64").Value = 3441.7856
$ws.Range("L64").Value = 4833.3335
$ws.Range("M64").Value = -3193.7856
$ws.Range("N64").Value = -5329.3335
$ws.Range("H67").Value = 3859.25
$ws.Range("I67").Value = 3441.7856
$ws.Range("J67").Value = 4833.3335
$ws.Range("K67").Value = 3441.7856
$ws.Range("L67").Value = 4833.3335
$ws.Range("M67").Value = -2583.7856
$ws.Range("N67").Value = -6549.3335
$ws.Range("H74").Value = 5109.9
$ws.Range("I74").Value = 4675
$ws.Range("K74").Value = 4675
$ws.Range("M74").Value = -3739
$ws.Range("H77").Value = 5109.9
$ws.Range("I77").Value = 4675
$ws.Range("K77").Value = 23375
$ws.Range("M77").Value = -18695
$ws.Range("H113").Value = 1833.3334
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1833.3334
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1420.6666
$ws.Range("H135").Value = 150003540
$ws.Range("I135").Value = 83337230
$ws.Range("K135").Value = 750035070
$ws.Range("M135").Value = -750032535
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1434.6177
$ws.Range("I2").Value = 1432.6364
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 1432.6364
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -1319.6364
$ws.Range("N2").Value = -1726
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H116").Value = 1434.6177
$ws.Range("I116").Value = 1432.6364
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 1432.6364
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 861.3635999999999
$ws.Range("N116").Value = -6088
$ws.Range("H122").Value = 20834584
$ws.Range("I122").Value = 1166.3334
$ws.Range("K122").Value = 3499.0002
$ws.Range("M122").Value = -1049.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1434.6177
$ws.Range("I3").Value = 1432.6364
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1432.6364
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -1318.6364
$ws.Range("N3").Value = -1728
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0
$ws.Range("H86").Value = 1848.8718
$ws.Range("I86").Value = 1612.303
$ws.Range("J86").Value = 3150
$ws.Range("K86").Value = 1612.303
$ws.Range("L86").Value = 3150
$ws.Range("M86").Value = -489.3030000000001
$ws.Range("N86").Value = -5396
$ws.Range("H89").Value = 1848.8718
$ws.Range("I89").Value = 1612.303
$ws.Range("J89").Value = 3150
$ws.Range("K89").Value = 8061.515
$ws.Range("L89").Value = 15750
$ws.Range("M89").Value = -2445.515
$ws.Range("N89").Value = -26982
$ws.Range("H134").Value = 2957.75
$ws.Range("I134").Value = 2626.2666
$ws.Range("J134").Value = 3340.2307
$ws.Range("K134").Value = 7878.7998
$ws.Range("L134").Value = 10020.6921
$ws.Range("M134").Value = -5343.7998
$ws.Range("N134").Value = -15090.6921
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 349.5
$ws.Range("I22").Value = 299.33334
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 299.33334
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 50.66665999999998
$ws.Range("N22").Value = -1200
$ws.Range("H99").Value = 2381.818
$ws.Range("I99").Value = 1766.6666
$ws.Range("J99").Value = 3120
$ws.Range("K99").Value = 1766.6666
$ws.Range("L99").Value = 3120
$ws.Range("M99").Value = -268.6666
$ws.Range("N99").Value = -6116
$ws.Range("H126").Value = 2381.818
$ws.Range("I126").Value = 1766.6666
$ws.Range("J126").Value = 3120
$ws.Range("K126").Value = 5299.9998
$ws.Range("L126").Value = 9360
$ws.Range("M126").Value = -2829.9998
$ws.Range("N126").Value = -14300
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 175.9
$ws.Range("I11").Value = 173.42857
$ws.Range("J11").Value = 181.66667
$ws.Range("K11").Value = 520.28571
$ws.Range("L11").Value = 545.00001
$ws.Range("M11").Value = -380.28571
$ws.Range("N11").Value = -825.00001
$ws.Range("H131").Value = 15049.164
$ws.Range("I131").Value = 996.9375
$ws.Range("J131").Value = 19457.705
$ws.Range("K131").Value = 2990.8125
$ws.Range("L131").Value = 58373.11500000001
$ws.Range("M131").Value = 2049.1875
$ws.Range("N131").Value = -68453.11500000001
$ws.Range("H134").Value = 2592.675
$ws.Range("I134").Value = 1745.069
$ws.Range("J134").Value = 4827.273
$ws.Range("K134").Value = 5235.207
$ws.Range("L134").Value = 14481.819
$ws.Range("M134").Value = -165.2070000000003
$ws.Range("N134").Value = -24621.819
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1386.5714
$ws.Range("I13").Value = 1386.5714
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1386.5714
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1247.5714
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("H80").Value = 7169.615
$ws.Range("I80").Value = 3335
$ws.Range("J80").Value = 8320
$ws.Range("K80").Value = 3335
$ws.Range("L80").Value = 8320
$ws.Range("M80").Value = -2337
$ws.Range("N80").Value = -10316
$ws.Range("H83").Value = 7169.615
$ws.Range("I83").Value = 3335
$ws.Range("J83").Value = 8320
$ws.Range("K83").Value = 16675
$ws.Range("L83").Value = 41600
$ws.Range("M83").Value = -11683
$ws.Range("N83").Value = -51584
$ws.Range("H122").Value = 5915.6
$ws.Range("I122").Value = 10447.5
$ws.Range("J122").Value = 2894.3333
$ws.Range("K122").Value = 31342.5
$ws.Range("L122").Value = 8682.999899999999
$ws.Range("M122").Value = -28892.5
$ws.Range("N122").Value = -13582.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 277.16666
$ws.Range("I22").Value = 249.5
$ws.Range("J22").Value = 291
$ws.Range("K22").Value = 249.5
$ws.Range("L22").Value = 291
$ws.Range("M22").Value = 45.5
$ws.Range("N22").Value = -881
$ws.Range("H24").Value = 25000
$ws.Range("J24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("N24").Value = -25686
$ws.Range("H27").Value = 277.16666
$ws.Range("I27").Value = 249.5
$ws.Range("J27").Value = 291
$ws.Range("K27").Value = 249.5
$ws.Range("L27").Value = 291
$ws.Range("M27").Value = -142.5
$ws.Range("N27").Value = -505
$ws.Range("I40").Value = 2783.1667
$ws.Range("J40").Value = 4120.75
$ws.Range("K40").Value = 2783.1667
$ws.Range("L40").Value = 4120.75
$ws.Range("M40").Value = -2647.1667
$ws.Range("N40").Value = -4392.75
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H82").Value = 1979.6
$ws.Range("I82").Value = 1888.4445
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 1888.4445
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -1527.4445
$ws.Range("N82").Value = -3522
$ws.Range("H85").Value = 1979.6
$ws.Range("I85").Value = 1888.4445
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 1888.4445
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = -640.4445000000001
$ws.Range("N85").Value = -5296
$ws.Range("H141").Value = 21194.834
$ws.Range("J141").Value = 21194.834
$ws.Range("L141").Value = 21194.834
$ws.Range("N141").Value = -31554.834
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10529492
$ws.Range("I81").Value = 1745
$ws.Range("K81").Value = 3490
$ws.Range("M81").Value = -2429
$ws.Range("H84").Value = 10529492
$ws.Range("I84").Value = 1745
$ws.Range("K84").Value = 17450
$ws.Range("M84").Value = -12146
$ws.Range("H113").Value = 749.5517
$ws.Range("I113").Value = 418.2143
$ws.Range("J113").Value = 1058.8
$ws.Range("K113").Value = 1254.6429
$ws.Range("L113").Value = 3176.4
$ws.Range("M113").Value = 915.3571000000002
$ws.Range("N113").Value = -7516.4
